# Atualizado por script em 31-10-2023 15:01
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match order for the Aug 20/21 2023 double round (rows 21 & 22) ---
# Columns F:V hold the per-match data (B:E are shared/identical for the pair).
$tmpF = $ws.Cells.Item(21,6).Value2
$tmpG = $ws.Cells.Item(21,7).Value2
$tmpH = $ws.Cells.Item(21,8).Value2
$tmpI = $ws.Cells.Item(21,9).Value2
$tmpJ = $ws.Cells.Item(21,10).Value2
$tmpK = $ws.Cells.Item(21,11).Value2
$tmpL = $ws.Cells.Item(21,12).Value2
$tmpM = $ws.Cells.Item(21,13).Value2
$tmpN = $ws.Cells.Item(21,14).Value2
$tmpO = $ws.Cells.Item(21,15).Value2
$tmpP = $ws.Cells.Item(21,16).Value2
$tmpQ = $ws.Cells.Item(21,17).Value2
$tmpR = $ws.Cells.Item(21,18).Value2
$tmpS = $ws.Cells.Item(21,19).Value2
$tmpT = $ws.Cells.Item(21,20).Value2
$tmpU = $ws.Cells.Item(21,21).Value2
$tmpV = $ws.Cells.Item(21,22).Value2

$ws.Cells.Item(21,6).Value = $ws.Cells.Item(22,6).Value2
$ws.Cells.Item(21,7).Value = $ws.Cells.Item(22,7).Value2
$ws.Cells.Item(21,8).Value = $ws.Cells.Item(22,8).Value2
$ws.Cells.Item(21,9).Value = $ws.Cells.Item(22,9).Value2
$ws.Cells.Item(21,10).Value = $ws.Cells.Item(22,10).Value2
$ws.Cells.Item(21,11).Value = $ws.Cells.Item(22,11).Value2
$ws.Cells.Item(21,12).Value = $ws.Cells.Item(22,12).Value2
$ws.Cells.Item(21,13).Value = $ws.Cells.Item(22,13).Value2
$ws.Cells.Item(21,14).Value = $ws.Cells.Item(22,14).Value2
$ws.Cells.Item(21,15).Value = $ws.Cells.Item(22,15).Value2
$ws.Cells.Item(21,16).Value = $ws.Cells.Item(22,16).Value2
$ws.Cells.Item(21,17).Value = $ws.Cells.Item(22,17).Value2
$ws.Cells.Item(21,18).Value = $ws.Cells.Item(22,18).Value2
$ws.Cells.Item(21,19).Value = $ws.Cells.Item(22,19).Value2
$ws.Cells.Item(21,20).Value = $ws.Cells.Item(22,20).Value2
$ws.Cells.Item(21,21).Value = $ws.Cells.Item(22,21).Value2
$ws.Cells.Item(21,22).Value = $ws.Cells.Item(22,22).Value2

$ws.Cells.Item(22,6).Value = $tmpF
$ws.Cells.Item(22,7).Value = $tmpG
$ws.Cells.Item(22,8).Value = $tmpH
$ws.Cells.Item(22,9).Value = $tmpI
$ws.Cells.Item(22,10).Value = $tmpJ
$ws.Cells.Item(22,11).Value = $tmpK
$ws.Cells.Item(22,12).Value = $tmpL
$ws.Cells.Item(22,13).Value = $tmpM
$ws.Cells.Item(22,14).Value = $tmpN
$ws.Cells.Item(22,15).Value = $tmpO
$ws.Cells.Item(22,16).Value = $tmpP
$ws.Cells.Item(22,17).Value = $tmpQ
$ws.Cells.Item(22,18).Value = $tmpR
$ws.Cells.Item(22,19).Value = $tmpS
$ws.Cells.Item(22,20).Value = $tmpT
$ws.Cells.Item(22,21).Value = $tmpU
$ws.Cells.Item(22,22).Value = $tmpV

# --- Swap match order for the Oct 16/2023 double round (rows 79 & 80) ---
$tmpF = $ws.Cells.Item(79,6).Value2
$tmpG = $ws.Cells.Item(79,7).Value2
$tmpH = $ws.Cells.Item(79,8).Value2
$tmpI = $ws.Cells.Item(79,9).Value2
$tmpJ = $ws.Cells.Item(79,10).Value2
$tmpK = $ws.Cells.Item(79,11).Value2
$tmpL = $ws.Cells.Item(79,12).Value2
$tmpM = $ws.Cells.Item(79,13).Value2
$tmpN = $ws.Cells.Item(79,14).Value2
$tmpO = $ws.Cells.Item(79,15).Value2
$tmpP = $ws.Cells.Item(79,16).Value2
$tmpQ = $ws.Cells.Item(79,17).Value2
$tmpR = $ws.Cells.Item(79,18).Value2
$tmpS = $ws.Cells.Item(79,19).Value2
$tmpT = $ws.Cells.Item(79,20).Value2
$tmpU = $ws.Cells.Item(79,21).Value2
$tmpV = $ws.Cells.Item(79,22).Value2

$ws.Cells.Item(79,6).Value = $ws.Cells.Item(80,6).Value2
$ws.Cells.Item(79,7).Value = $ws.Cells.Item(80,7).Value2
$ws.Cells.Item(79,8).Value = $ws.Cells.Item(80,8).Value2
$ws.Cells.Item(79,9).Value = $ws.Cells.Item(80,9).Value2
$ws.Cells.Item(79,10).Value = $ws.Cells.Item(80,10).Value2
$ws.Cells.Item(79,11).Value = $ws.Cells.Item(80,11).Value2
$ws.Cells.Item(79,12).Value = $ws.Cells.Item(80,12).Value2
$ws.Cells.Item(79,13).Value = $ws.Cells.Item(80,13).Value2
$ws.Cells.Item(79,14).Value = $ws.Cells.Item(80,14).Value2
$ws.Cells.Item(79,15).Value = $ws.Cells.Item(80,15).Value2
$ws.Cells.Item(79,16).Value = $ws.Cells.Item(80,16).Value2
$ws.Cells.Item(79,17).Value = $ws.Cells.Item(80,17).Value2
$ws.Cells.Item(79,18).Value = $ws.Cells.Item(80,18).Value2
$ws.Cells.Item(79,19).Value = $ws.Cells.Item(80,19).Value2
$ws.Cells.Item(79,20).Value = $ws.Cells.Item(80,20).Value2
$ws.Cells.Item(79,21).Value = $ws.Cells.Item(80,21).Value2
$ws.Cells.Item(79,22).Value = $ws.Cells.Item(80,22).Value2

$ws.Cells.Item(80,6).Value = $tmpF
$ws.Cells.Item(80,7).Value = $tmpG
$ws.Cells.Item(80,8).Value = $tmpH
$ws.Cells.Item(80,9).Value = $tmpI
$ws.Cells.Item(80,10).Value = $tmpJ
$ws.Cells.Item(80,11).Value = $tmpK
$ws.Cells.Item(80,12).Value = $tmpL
$ws.Cells.Item(80,13).Value = $tmpM
$ws.Cells.Item(80,14).Value = $tmpN
$ws.Cells.Item(80,15).Value = $tmpO
$ws.Cells.Item(80,16).Value = $tmpP
$ws.Cells.Item(80,17).Value = $tmpQ
$ws.Cells.Item(80,18).Value = $tmpR
$ws.Cells.Item(80,19).Value = $tmpS
$ws.Cells.Item(80,20).Value = $tmpT
$ws.Cells.Item(80,21).Value = $tmpU
$ws.Cells.Item(80,22).Value = $tmpV

# --- Append 9 new match rows (83-91), matching the formatting of the last
#     existing data row (82) for columns A (index) and E (datetime) ---
$ws.Range("A82:V82").Copy()
$ws.Range("A83:V91").PasteSpecial(-4122)

$ws.Cells.Item(83,1).Value = 82
$ws.Cells.Item(83,2).Value = "turkey"
$ws.Cells.Item(83,3).Value = "1-lig"
$ws.Cells.Item(83,4).Value = "2023-2024"
$ws.Cells.Item(83,5).Value = 45226.79166666666
$ws.Cells.Item(83,6).Value = "Sanliurfaspor"
$ws.Cells.Item(83,7).Value = 0
$ws.Cells.Item(83,8).Value = "Giresunspor"
$ws.Cells.Item(83,9).Value = 0
$ws.Cells.Item(83,10).Value = 1.76
$ws.Cells.Item(83,11).Value = "26/10/2023 13:55"
$ws.Cells.Item(83,12).Value = 1.64
$ws.Cells.Item(83,13).Value = "27/10/2023 18:59"
$ws.Cells.Item(83,14).Value = 3.65
$ws.Cells.Item(83,15).Value = "26/10/2023 13:55"
$ws.Cells.Item(83,16).Value = 3.97
$ws.Cells.Item(83,17).Value = "27/10/2023 18:59"
$ws.Cells.Item(83,18).Value = 4.24
$ws.Cells.Item(83,19).Value = "26/10/2023 13:55"
$ws.Cells.Item(83,20).Value = 5.33
$ws.Cells.Item(83,21).Value = "27/10/2023 18:59"
$ws.Cells.Item(83,22).Value = "https://www.betexplorer.com/football/turkey/1-lig/sanliurfaspor-giresunspor/pSpLPABU/"
$ws.Cells.Item(84,1).Value = 83
$ws.Cells.Item(84,2).Value = "turkey"
$ws.Cells.Item(84,3).Value = "1-lig"
$ws.Cells.Item(84,4).Value = "2023-2024"
$ws.Cells.Item(84,5).Value = 45227.52083333334
$ws.Cells.Item(84,6).Value = "Keciorengucu"
$ws.Cells.Item(84,7).Value = 0
$ws.Cells.Item(84,8).Value = "Goztepe"
$ws.Cells.Item(84,9).Value = 1
$ws.Cells.Item(84,10).Value = 2.56
$ws.Cells.Item(84,11).Value = "21/10/2023 18:13"
$ws.Cells.Item(84,12).Value = 3.25
$ws.Cells.Item(84,13).Value = "28/10/2023 12:21"
$ws.Cells.Item(84,14).Value = 3.18
$ws.Cells.Item(84,15).Value = "21/10/2023 18:13"
$ws.Cells.Item(84,16).Value = 3.11
$ws.Cells.Item(84,17).Value = "28/10/2023 12:21"
$ws.Cells.Item(84,18).Value = 2.87
$ws.Cells.Item(84,19).Value = "21/10/2023 18:13"
$ws.Cells.Item(84,20).Value = 2.41
$ws.Cells.Item(84,21).Value = "28/10/2023 12:21"
$ws.Cells.Item(84,22).Value = "https://www.betexplorer.com/football/turkey/1-lig/keciorengucu-goztepe/ziQUqEdn/"
$ws.Cells.Item(85,1).Value = 84
$ws.Cells.Item(85,2).Value = "turkey"
$ws.Cells.Item(85,3).Value = "1-lig"
$ws.Cells.Item(85,4).Value = "2023-2024"
$ws.Cells.Item(85,5).Value = 45227.52083333334
$ws.Cells.Item(85,6).Value = "Tuzlaspor"
$ws.Cells.Item(85,7).Value = 0
$ws.Cells.Item(85,8).Value = "Erzurumspor"
$ws.Cells.Item(85,9).Value = 1
$ws.Cells.Item(85,10).Value = 2.42
$ws.Cells.Item(85,11).Value = "21/10/2023 12:42"
$ws.Cells.Item(85,12).Value = 2.46
$ws.Cells.Item(85,13).Value = "28/10/2023 12:28"
$ws.Cells.Item(85,14).Value = 3.07
$ws.Cells.Item(85,15).Value = "21/10/2023 12:42"
$ws.Cells.Item(85,16).Value = 3.16
$ws.Cells.Item(85,17).Value = "28/10/2023 12:28"
$ws.Cells.Item(85,18).Value = 3.01
$ws.Cells.Item(85,19).Value = "21/10/2023 12:42"
$ws.Cells.Item(85,20).Value = 3.12
$ws.Cells.Item(85,21).Value = "28/10/2023 12:28"
$ws.Cells.Item(85,22).Value = "https://www.betexplorer.com/football/turkey/1-lig/tuzlaspor-erzurumspor-fk/d6EwrhRb/"
$ws.Cells.Item(86,1).Value = 85
$ws.Cells.Item(86,2).Value = "turkey"
$ws.Cells.Item(86,3).Value = "1-lig"
$ws.Cells.Item(86,4).Value = "2023-2024"
$ws.Cells.Item(86,5).Value = 45227.625
$ws.Cells.Item(86,6).Value = "Boluspor"
$ws.Cells.Item(86,7).Value = 2
$ws.Cells.Item(86,8).Value = "Bodrumspor"
$ws.Cells.Item(86,9).Value = 1
$ws.Cells.Item(86,10).Value = 3.43
$ws.Cells.Item(86,11).Value = "22/10/2023 15:12"
$ws.Cells.Item(86,12).Value = 3.56
$ws.Cells.Item(86,13).Value = "28/10/2023 14:54"
$ws.Cells.Item(86,14).Value = 3.29
$ws.Cells.Item(86,15).Value = "22/10/2023 15:12"
$ws.Cells.Item(86,16).Value = 3.04
$ws.Cells.Item(86,17).Value = "28/10/2023 14:54"
$ws.Cells.Item(86,18).Value = 2.18
$ws.Cells.Item(86,19).Value = "22/10/2023 15:12"
$ws.Cells.Item(86,20).Value = 2.3
$ws.Cells.Item(86,21).Value = "28/10/2023 14:54"
$ws.Cells.Item(86,22).Value = "https://www.betexplorer.com/football/turkey/1-lig/boluspor-bodrumspor/MmMQpfst/"
$ws.Cells.Item(87,1).Value = 86
$ws.Cells.Item(87,2).Value = "turkey"
$ws.Cells.Item(87,3).Value = "1-lig"
$ws.Cells.Item(87,4).Value = "2023-2024"
$ws.Cells.Item(87,5).Value = 45227.75
$ws.Cells.Item(87,6).Value = "Kocaelispor"
$ws.Cells.Item(87,7).Value = 2
$ws.Cells.Item(87,8).Value = "Genclerbirligi"
$ws.Cells.Item(87,9).Value = 2
$ws.Cells.Item(87,10).Value = 2.13
$ws.Cells.Item(87,11).Value = "26/10/2023 15:12"
$ws.Cells.Item(87,12).Value = 1.87
$ws.Cells.Item(87,13).Value = "28/10/2023 17:58"
$ws.Cells.Item(87,14).Value = 3.42
$ws.Cells.Item(87,15).Value = "26/10/2023 15:12"
$ws.Cells.Item(87,16).Value = 3.61
$ws.Cells.Item(87,17).Value = "28/10/2023 17:58"
$ws.Cells.Item(87,18).Value = 3.42
$ws.Cells.Item(87,19).Value = "26/10/2023 15:12"
$ws.Cells.Item(87,20).Value = 4.24
$ws.Cells.Item(87,21).Value = "28/10/2023 17:58"
$ws.Cells.Item(87,22).Value = "https://www.betexplorer.com/football/turkey/1-lig/kocaelispor-genclerbirligi/IFFZrYBh/"
$ws.Cells.Item(88,1).Value = 87
$ws.Cells.Item(88,2).Value = "turkey"
$ws.Cells.Item(88,3).Value = "1-lig"
$ws.Cells.Item(88,4).Value = "2023-2024"
$ws.Cells.Item(88,5).Value = 45228.47916666666
$ws.Cells.Item(88,6).Value = "Manisa FK"
$ws.Cells.Item(88,7).Value = 2
$ws.Cells.Item(88,8).Value = "Corum"
$ws.Cells.Item(88,9).Value = 2
$ws.Cells.Item(88,10).Value = 1.93
$ws.Cells.Item(88,11).Value = "26/10/2023 15:12"
$ws.Cells.Item(88,12).Value = 1.87
$ws.Cells.Item(88,13).Value = "29/10/2023 11:02"
$ws.Cells.Item(88,14).Value = 3.54
$ws.Cells.Item(88,15).Value = "26/10/2023 15:12"
$ws.Cells.Item(88,16).Value = 3.56
$ws.Cells.Item(88,17).Value = "29/10/2023 11:02"
$ws.Cells.Item(88,18).Value = 3.9
$ws.Cells.Item(88,19).Value = "26/10/2023 15:12"
$ws.Cells.Item(88,20).Value = 4.34
$ws.Cells.Item(88,21).Value = "29/10/2023 11:01"
$ws.Cells.Item(88,22).Value = "https://www.betexplorer.com/football/turkey/1-lig/manisa-fk-corum-fk/ne7jujBH/"
$ws.Cells.Item(89,1).Value = 88
$ws.Cells.Item(89,2).Value = "turkey"
$ws.Cells.Item(89,3).Value = "1-lig"
$ws.Cells.Item(89,4).Value = "2023-2024"
$ws.Cells.Item(89,5).Value = 45228.58333333334
$ws.Cells.Item(89,6).Value = "Adanaspor AS"
$ws.Cells.Item(89,7).Value = 1
$ws.Cells.Item(89,8).Value = "Umraniyespor"
$ws.Cells.Item(89,9).Value = 0
$ws.Cells.Item(89,10).Value = 2.55
$ws.Cells.Item(89,11).Value = "26/10/2023 15:12"
$ws.Cells.Item(89,12).Value = 3.46
$ws.Cells.Item(89,13).Value = "29/10/2023 13:52"
$ws.Cells.Item(89,14).Value = 3.24
$ws.Cells.Item(89,15).Value = "26/10/2023 15:12"
$ws.Cells.Item(89,16).Value = 3.39
$ws.Cells.Item(89,17).Value = "29/10/2023 13:52"
$ws.Cells.Item(89,18).Value = 2.83
$ws.Cells.Item(89,19).Value = "26/10/2023 15:12"
$ws.Cells.Item(89,20).Value = 2.17
$ws.Cells.Item(89,21).Value = "29/10/2023 13:52"
$ws.Cells.Item(89,22).Value = "https://www.betexplorer.com/football/turkey/1-lig/adanaspor-as-umraniyespor/6ynHQjdO/"
$ws.Cells.Item(90,1).Value = 89
$ws.Cells.Item(90,2).Value = "turkey"
$ws.Cells.Item(90,3).Value = "1-lig"
$ws.Cells.Item(90,4).Value = "2023-2024"
$ws.Cells.Item(90,5).Value = 45228.58333333334
$ws.Cells.Item(90,6).Value = "Altay"
$ws.Cells.Item(90,7).Value = 1
$ws.Cells.Item(90,8).Value = "Sakaryaspor"
$ws.Cells.Item(90,9).Value = 2
$ws.Cells.Item(90,10).Value = 2.9
$ws.Cells.Item(90,11).Value = "23/10/2023 19:12"
$ws.Cells.Item(90,12).Value = 4.56
$ws.Cells.Item(90,13).Value = "29/10/2023 13:52"
$ws.Cells.Item(90,14).Value = 3.32
$ws.Cells.Item(90,15).Value = "23/10/2023 19:12"
$ws.Cells.Item(90,16).Value = 3.66
$ws.Cells.Item(90,17).Value = "29/10/2023 13:52"
$ws.Cells.Item(90,18).Value = 2.45
$ws.Cells.Item(90,19).Value = "23/10/2023 19:12"
$ws.Cells.Item(90,20).Value = 1.8
$ws.Cells.Item(90,21).Value = "29/10/2023 13:52"
$ws.Cells.Item(90,22).Value = "https://www.betexplorer.com/football/turkey/1-lig/altay-sakaryaspor/4WHotWeB/"
$ws.Cells.Item(91,1).Value = 90
$ws.Cells.Item(91,2).Value = "turkey"
$ws.Cells.Item(91,3).Value = "1-lig"
$ws.Cells.Item(91,4).Value = "2023-2024"
$ws.Cells.Item(91,5).Value = 45228.70833333334
$ws.Cells.Item(91,6).Value = "Eyupspor"
$ws.Cells.Item(91,7).Value = 0
$ws.Cells.Item(91,8).Value = "Bandirmaspor"
$ws.Cells.Item(91,9).Value = 1
$ws.Cells.Item(91,10).Value = 1.61
$ws.Cells.Item(91,11).Value = "22/10/2023 20:15"
$ws.Cells.Item(91,12).Value = 1.65
$ws.Cells.Item(91,13).Value = "29/10/2023 16:55"
$ws.Cells.Item(91,14).Value = 3.94
$ws.Cells.Item(91,15).Value = "22/10/2023 20:15"
$ws.Cells.Item(91,16).Value = 3.81
$ws.Cells.Item(91,17).Value = "29/10/2023 16:55"
$ws.Cells.Item(91,18).Value = 5.41
$ws.Cells.Item(91,19).Value = "22/10/2023 20:15"
$ws.Cells.Item(91,20).Value = 5.53
$ws.Cells.Item(91,21).Value = "29/10/2023 16:55"
$ws.Cells.Item(91,22).Value = "https://www.betexplorer.com/football/turkey/1-lig/eyupspor-bandirmaspor/QuJssCt5/"
